# Insert a new data row at row 342 (pushes existing rows 342-385 down to 343-386)
# and populate it with a new "Papa" price observation for
# "Vega Monumental Concepción" / Bíobío.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(342).Insert()

$ws.Range("A342").Value = 11
$ws.Range("B342").Value = "Vega Monumental Concepción"
$ws.Range("C342").Value = "Bíobío"
$ws.Range("D342").Value = 44946
$ws.Range("E342").Value = 8
$ws.Range("F342").Value = 100114001
$ws.Range("G342").Value = "Papa"
$ws.Range("H342").Value = "Asterix"
$ws.Range("I342").Value = "1a nueva(o)"
$ws.Range("J342").Value = 220
$ws.Range("K342").Value = 10000
$ws.Range("L342").Value = 11000
$ws.Range("M342").Value = 10545
$ws.Range("N342").Value = "`$/saco 25 kilos"
$ws.Range("O342").Value = "Región de La Araucanía"
$ws.Range("P342").Value = 422
$ws.Range("Q342").Value = 25
$ws.Range("R342").Value = "Hortaliza"
